$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 gets the values previously in row 3
$ws.Range("D2").Value = 44421
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 24000
$ws.Range("O2").Value = 24000
$ws.Range("P2").Value = 24000
$ws.Range("S2").Value = 1200

# Row 3 gets the values previously in row 2
$ws.Range("D3").Value = 44291
$ws.Range("M3").Value = 15
$ws.Range("N3").Value = 23000
$ws.Range("O3").Value = 23000
$ws.Range("P3").Value = 23000
$ws.Range("S3").Value = 1150
